# Second Run for Stock Price Prediction ML project w Basic RF
# Replace the 14-column (A:N) results table with a new 9-column (A:I) table:
# the old Train/Test/Val metric columns are dropped in favor of a
# Test/Val-only layout, and all metric values are refreshed with the
# second-run numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old RMSE(Train)/DA(Train) etc. columns (J:N) entirely so the
# sheet's dimension shrinks back down to A1:I9.
$ws.Range("J1:N9").Clear()

# Store the refreshed metrics as literal text (matching how the rest of
# the sheet already stores its numbers) so values like "0.3000" or
# "214506.3424" round-trip exactly instead of being reinterpreted as
# floating point numbers.
$ws.Range("B2:I9").NumberFormat = "@"

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "stock"
$ws.Range("B1").Value = "R2 (Test)"
$ws.Range("C1").Value = "R2 (Val)"
$ws.Range("D1").Value = "MSE (Test)"
$ws.Range("E1").Value = "MSE (Val)"
$ws.Range("F1").Value = "RMSE (Test)"
$ws.Range("G1").Value = "RMSE (Val)"
$ws.Range("H1").Value = "DA (Test)"
$ws.Range("I1").Value = "DA (Val)"

# --- Row 2: AMZN --------------------------------------------------------
$ws.Range("B2").Value = "-14.8653"
$ws.Range("C2").Value = "-0.0123"
$ws.Range("D2").Value = "7020.8198"
$ws.Range("E2").Value = "897.9986"
$ws.Range("F2").Value = "83.7903"
$ws.Range("G2").Value = "29.9666"
$ws.Range("H2").Value = "0.3568"
$ws.Range("I2").Value = "0.4548"

# --- Row 3: AAPL --------------------------------------------------------
$ws.Range("B3").Value = "-14.1246"
$ws.Range("C3").Value = "-0.8470"
$ws.Range("D3").Value = "5669.5616"
$ws.Range("E3").Value = "635.8958"
$ws.Range("F3").Value = "75.2965"
$ws.Range("G3").Value = "25.2170"
$ws.Range("H3").Value = "0.3324"
$ws.Range("I3").Value = "0.5824"

# --- Row 4: GOOG --------------------------------------------------------
$ws.Range("B4").Value = "-5.7311"
$ws.Range("C4").Value = "-0.1985"
$ws.Range("D4").Value = "6602.9353"
$ws.Range("E4").Value = "636.0069"
$ws.Range("F4").Value = "81.2584"
$ws.Range("G4").Value = "25.2192"
$ws.Range("H4").Value = "0.3000"
$ws.Range("I4").Value = "0.3852"

# --- Row 5: NDAQ --------------------------------------------------------
$ws.Range("B5").Value = "-3.8267"
$ws.Range("C5").Value = "-0.2580"
$ws.Range("D5").Value = "478.5966"
$ws.Range("E5").Value = "28.2750"
$ws.Range("F5").Value = "21.8769"
$ws.Range("G5").Value = "5.3174"
$ws.Range("H5").Value = "0.3838"
$ws.Range("I5").Value = "0.6682"

# --- Row 6: META --------------------------------------------------------
$ws.Range("B6").Value = "-25.3776"
$ws.Range("C6").Value = "-0.8796"
$ws.Range("D6").Value = "214506.3424"
$ws.Range("E6").Value = "28668.3409"
$ws.Range("F6").Value = "463.1483"
$ws.Range("G6").Value = "169.3173"
$ws.Range("H6").Value = "0.3351"
$ws.Range("I6").Value = "0.3248"

# --- Row 7: TSLA --------------------------------------------------------
$ws.Range("B7").Value = "-0.2557"
$ws.Range("C7").Value = "-3.0053"
$ws.Range("D7").Value = "8262.0997"
$ws.Range("E7").Value = "7188.2126"
$ws.Range("F7").Value = "90.8961"
$ws.Range("G7").Value = "84.7833"
$ws.Range("H7").Value = "0.6297"
$ws.Range("I7").Value = "0.4640"

# --- Row 8: INTC --------------------------------------------------------
$ws.Range("B8").Value = "-5.7259"
$ws.Range("C8").Value = "-0.2654"
$ws.Range("D8").Value = "217.3009"
$ws.Range("E8").Value = "53.0084"
$ws.Range("F8").Value = "14.7411"
$ws.Range("G8").Value = "7.2807"
$ws.Range("H8").Value = "0.5919"
$ws.Range("I8").Value = "0.6752"

# --- Row 9: AMD ---------------------------------------------------------
$ws.Range("B9").Value = "-0.8899"
$ws.Range("C9").Value = "0.0426"
$ws.Range("D9").Value = "2612.9322"
$ws.Range("E9").Value = "1376.8297"
$ws.Range("F9").Value = "51.1168"
$ws.Range("G9").Value = "37.1057"
$ws.Range("H9").Value = "0.5595"
$ws.Range("I9").Value = "0.5940"
